{"js": "// Add a new \"Ref: SYS_155\" paragraph right after the \"one for jumping\"\n// bullet item (inside the \"Requirements\" section), matching the\n// commit \"add undefined referenced requirement in docx\".\n//\n// The new paragraph uses the base/\"Normal\" paragraph style (style0) and\n// carries no list numbering, so rather than inserting immediately after\n// the bulleted \"one for jumping\" paragraph (which would make the new\n// paragraph inherit its bullet numbering), we insert the new paragraph\n// right before the paragraph that already follows it \u2014 that paragraph\n// is a plain \"Normal\"-style paragraph, so the new one naturally picks\n// up the same (numbering-free) formatting.\n\nconst body = context.document.body;\n\n// Locate the \"one for jumping\" paragraph robustly via search instead of\n// a hard-coded paragraph index.\nconst searchResults = body.search(\"one for jumping\", { matchCase: true });\nsearchResults.load(\"items\");\nawait context.sync();\n\nconst anchorParagraph = searchResults.items[0].paragraphs.getFirst();\nawait context.sync();\n\nconst followingParagraph = anchorParagraph.getNext();\nfollowingParagraph.insertParagraph(\"Ref: SYS_155\", Word.InsertLocation.before);\n\nawait context.sync();\n", "ps1": "# Add a new \"Ref: SYS_155\" paragraph right after the \"one for jumping\"\n# bullet item (inside the \"Requirements\" section), matching the\n# commit \"add undefined referenced requirement in docx\".\n#\n# The new paragraph uses the base/\"Normal\" paragraph style (style0) and\n# carries no list numbering, so rather than inserting immediately after\n# the bulleted \"one for jumping\" paragraph (which would make the new\n# paragraph inherit its bullet numbering), we insert the new paragraph\n# right before the paragraph that already follows it \u2014 that paragraph\n# is a plain \"Normal\"-style paragraph, so the new one naturally picks\n# up the same (numbering-free) formatting.\n\n$d = $word.ActiveDocument\n\n# Locate the \"one for jumping\" paragraph robustly via Find instead of a\n# hard-coded paragraph index.\n$searchRange = $d.Content\n$find = $searchRange.Find\n$find.ClearFormatting()\n$find.Text = \"one for jumping\"\n$found = $find.Execute()\n\n$allParagraphs = $d.Paragraphs\n$anchorIndex = -1\nfor ($i = 1; $i -le $allParagraphs.Count; $i++) {\n    $candidate = $allParagraphs.Item($i)\n    if ($candidate.Range.Start -le $searchRange.Start -and $candidate.Range.End -ge $searchRange.End) {\n        $anchorIndex = $i\n        break\n    }\n}\n\n$followingParagraph = $allParagraphs.Item($anchorIndex + 1)\n$insertionPoint = $followingParagraph.Range\n$insertionPoint.Collapse(1)  # wdCollapseStart\n$insertionPoint.InsertBefore(\"Ref: SYS_155`r\")\n"}
